$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Long POM reference text for "Build a executable JAR" (maven-assembly-plugin) ---
$pom1 = @'
<project xmlns="http://maven.apache.org/POM/4.0.0" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance"
  xsi:schemaLocation="http://maven.apache.org/POM/4.0.0 http://maven.apache.org/xsd/maven-4.0.0.xsd">
  <modelVersion>4.0.0</modelVersion>
  <groupId>com.myles.util.ping</groupId>
  <artifactId>LocationTeller</artifactId>
  <version>1.0-SNAPSHOT</version>
  <packaging>jar</packaging>
  <name>LocationTeller</name>
  <url>http://maven.apache.org</url>
  <properties>
    <project.build.sourceEncoding>UTF-8</project.build.sourceEncoding>
  </properties>
  <dependencies>
    <dependency>
      <groupId>junit</groupId>
      <artifactId>junit</artifactId>
      <version>3.8.1</version>
      <scope>test</scope>
    </dependency>
  </dependencies>
  <build>
    <plugins>
      <plugin>
        <artifactId>maven-assembly-plugin</artifactId>
        <configuration>
          <archive>
            <manifest>
              <mainClass>com.myles.util.ping.App</mainClass>
            </manifest>
          </archive>
          <descriptorRefs>
            <descriptorRef>jar-with-dependencies</descriptorRef>
          </descriptorRefs>
        </configuration>
      </plugin>
    </plugins>
  </build>
</project>

'@

# --- Long POM reference text for "Build a window executable (*.exe)" (launch4j) ---
$pom2 = @'
<project xmlns="http://maven.apache.org/POM/4.0.0" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance"
  xsi:schemaLocation="http://maven.apache.org/POM/4.0.0 http://maven.apache.org/xsd/maven-4.0.0.xsd">
  <modelVersion>4.0.0</modelVersion>
  <groupId>com.myles.util.ping</groupId>
  <artifactId>LocationTeller</artifactId>
  <version>1.0-SNAPSHOT</version>
  <packaging>jar</packaging>
  <name>LocationTeller</name>
  <url>http://maven.apache.org</url>
  <properties>
    <project.build.sourceEncoding>UTF-8</project.build.sourceEncoding>
  </properties>
  <dependencies>
    <dependency>
      <groupId>junit</groupId>
      <artifactId>junit</artifactId>
      <version>3.8.1</version>
      <scope>test</scope>
    </dependency>
  </dependencies>
  <build>
      <plugins>
          <plugin>
              <groupId>org.apache.maven.plugins</groupId>
              <artifactId>maven-compiler-plugin</artifactId>
              <version>2.5.1</version>
              <configuration>
                  <source>1.6</source>
                  <target>1.6</target>
              </configuration>
          </plugin>
          <plugin>
              <groupId>org.apache.maven.plugins</groupId>
              <artifactId>maven-shade-plugin</artifactId>
              <version>1.7.1</version>
              <executions>
                  <execution>
                      <phase>package</phase>
                      <goals>
                          <goal>shade</goal>
                      </goals>
                  </execution>
              </executions>
              <configuration>
                  <shadedArtifactAttached>true</shadedArtifactAttached>
                  <shadedClassifierName>shaded</shadedClassifierName>
                  <transformers>
                      <transformer implementation="org.apache.maven.plugins.shade.resource.ManifestResourceTransformer">
                          <mainClass>com.myles.util.ping.App</mainClass>
                      </transformer>
                  </transformers>
              </configuration>
          </plugin>
          <plugin>
              <groupId>com.akathist.maven.plugins.launch4j</groupId>
              <artifactId>launch4j-maven-plugin</artifactId>
              <version>1.5.1</version>
              <executions>
                  <execution>
                      <id>l4j-clui</id>
                      <phase>package</phase>
                      <goals>
                          <goal>launch4j</goal>
                      </goals>
                      <configuration>
                          <headerType>gui</headerType>
                          <jar>${project.build.directory}/${artifactId}-${version}-shaded.jar</jar>
                          <outfile>${project.build.directory}/myles.exe</outfile>
                          <downloadUrl>http://java.com/download</downloadUrl>
                          <classPath>
                              <mainClass>com.myles.util.ping.App</mainClass>
                              <preCp>anything</preCp>
                          </classPath>
                          <icon/>
                          <jre>
                              <minVersion>1.6.0</minVersion>
                              <jdkPreference>preferJre</jdkPreference>
                          </jre>
                      </configuration>
                  </execution>
              </executions>
          </plugin>
      </plugins>
  </build>
</project>

'@

$windowExeNote = "Build a window executable (*.exe)`n** not compete, need further work`n** use the launch4j plugin"

# --- Row 9: interactive mode archetype generation ---
$ws.Cells.Item(9, 1).Value = "Maven"
$ws.Cells.Item(9, 2).Value = "build a project with interactive mode"
$ws.Cells.Item(9, 3).Value = "mvn archetype:generate"
$ws.Range("B8:C8").Copy() | Out-Null
$ws.Range("B9:C9").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(9).RowHeight = 36

# --- Row 10: build a runnable/executable jar via assembly plugin ---
$ws.Cells.Item(10, 1).Value = "Maven"
$ws.Cells.Item(10, 2).Value = "Build a executable JAR"
$ws.Cells.Item(10, 3).Value = $pom1
$ws.Range("B8:C8").Copy() | Out-Null
$ws.Range("B10:C10").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(10).RowHeight = 36

# --- Row 11: build a windows executable via launch4j plugin ---
$ws.Cells.Item(11, 1).Value = "Maven"
$ws.Cells.Item(11, 2).Value = $windowExeNote
$ws.Cells.Item(11, 3).Value = $pom2
$ws.Range("B8:C8").Copy() | Out-Null
$ws.Range("B11:C11").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(11).RowHeight = 36

$excel.CutCopyMode = 0
$ws.Range("C10").Select() | Out-Null
